$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column: header in H1 (copy formatting from the existing
# header cell G1 so it matches the other column headers), value in H2.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
